$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("H106").Value = 3479.7
$ws.Range("I106").Value = 3644.2222
$ws.Range("J106").Value = 1999
$ws.Range("K106").Value = 3644.2222
$ws.Range("L106").Value = 1999
$ws.Range("M106").Value = -3013.2222
$ws.Range("N106").Value = -3261
$ws.Range("H116").Value = 1826.5
$ws.Range("J116").Value = 1153
$ws.Range("L116").Value = 1153
$ws.Range("N116").Value = -8037
$ws.Range("H133").Value = 67387.5
$ws.Range("J133").Value = 67387.5
$ws.Range("L133").Value = 67387.5
$ws.Range("N133").Value = -77507.5
$ws.Range("N105").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H57").Value = 6000
$ws.Range("I57").Value = 6000
$ws.Range("K57").Value = 6000
$ws.Range("M57").Value = -5516
$ws.Range("H126").Value = 12000
$ws.Range("I126").Value = 12000
$ws.Range("K126").Value = 36000
$ws.Range("M126").Value = -33530

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 266.33334
$ws.Range("I22").Value = 274.5
$ws.Range("K22").Value = 274.5
$ws.Range("M22").Value = -101.5
$ws.Range("H99").Value = 1212.8572
$ws.Range("I99").Value = 997.5
$ws.Range("J99").Value = 1500
$ws.Range("K99").Value = 997.5
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = 500.5
$ws.Range("N99").Value = -4496
$ws.Range("H113").Value = 12000
$ws.Range("I113").Value = 12000
$ws.Range("K113").Value = 12000
$ws.Range("M113").Value = -9830
$ws.Range("H128").Value = 5000
$ws.Range("I128").Value = 5000
$ws.Range("K128").Value = 15000
$ws.Range("M128").Value = -12510

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H76").Value = 11499.667
$ws.Range("I76").Value = 11499.667
$ws.Range("K76").Value = 11499.667
$ws.Range("M76").Value = -11184.667
$ws.Range("H79").Value = 11499.667
$ws.Range("I79").Value = 11499.667
$ws.Range("K79").Value = 11499.667
$ws.Range("M79").Value = -10407.667
$ws.Range("H99").Value = 1242.6666
$ws.Range("I99").Value = 1131.2
$ws.Range("K99").Value = 1131.2
$ws.Range("M99").Value = 366.8
$ws.Range("H126").Value = 1242.6666
$ws.Range("I126").Value = 1131.2
$ws.Range("K126").Value = 3393.6
$ws.Range("M126").Value = -923.6000000000004
$ws.Range("H134").Value = 2269.0476
$ws.Range("I134").Value = 2054.0881
$ws.Range("J134").Value = 3182.625
$ws.Range("K134").Value = 6162.2643
$ws.Range("L134").Value = 9547.875
$ws.Range("M134").Value = -3627.2643
$ws.Range("N134").Value = -14617.875

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3357.6924
$ws.Range("I55").Value = 1466.6666
$ws.Range("J55").Value = 3604.348
$ws.Range("K55").Value = 4399.9998
$ws.Range("L55").Value = 10813.044
$ws.Range("M55").Value = -4222.9998
$ws.Range("N55").Value = -11167.044
$ws.Range("H92").Value = 449
$ws.Range("I92").Value = 449
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1347
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -99
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("H95").Value = 6007.3335
$ws.Range("J95").Value = 5999
$ws.Range("L95").Value = 17997
$ws.Range("N95").Value = -22115
$ws.Range("H96").Value = 4525
$ws.Range("J96").Value = 4525
$ws.Range("L96").Value = 13575
$ws.Range("N96").Value = -17693
$ws.Range("H97").Value = 28876.75
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 28876.75
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 86630.25
$ws.Range("N97").Value = -87622.25
$ws.Range("H98").Value = 377.66666
$ws.Range("I98").Value = 326.14285
$ws.Range("K98").Value = 978.4285500000001
$ws.Range("M98").Value = 519.5714499999999
$ws.Range("H100").Value = 3985.6
$ws.Range("J100").Value = 3985.6
$ws.Range("L100").Value = 11956.8
$ws.Range("N100").Value = -13578.8
$ws.Range("H101").Value = 7227.25
$ws.Range("J101").Value = 7227.25
$ws.Range("L101").Value = 21681.75
$ws.Range("N101").Value = -26549.75
$ws.Range("H102").Value = 5365.9443
$ws.Range("J102").Value = 5408.4375
$ws.Range("L102").Value = 16225.3125
$ws.Range("N102").Value = -21093.3125
$ws.Range("H103").Value = 2161
$ws.Range("J103").Value = 2825.4
$ws.Range("L103").Value = 8476.200000000001
$ws.Range("N103").Value = -10234.2
$ws.Range("H104").Value = 480
$ws.Range("I104").Value = 480
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 1440
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = 1181
$ws.Range("H105").Value = 6318
$ws.Range("J105").Value = 6318
$ws.Range("L105").Value = 18954
$ws.Range("N105").Value = -24196
$ws.Range("H106").Value = 3721
$ws.Range("J106").Value = 3721
$ws.Range("L106").Value = 11163
$ws.Range("N106").Value = -13055
$ws.Range("H122").Value = 853.65216
$ws.Range("I122").Value = 426.33334
$ws.Range("K122").Value = 3837.00006
$ws.Range("M122").Value = -1387.00006
$ws.Range("H124").Value = 2920.3635
$ws.Range("I124").Value = 1453.3334
$ws.Range("J124").Value = 3470.5
$ws.Range("K124").Value = 4360.0002
$ws.Range("L124").Value = 10411.5
$ws.Range("M124").Value = 549.9997999999996
$ws.Range("N124").Value = -20231.5
$ws.Range("H125").Value = 3976.4614
$ws.Range("J125").Value = 4141.1665
$ws.Range("L125").Value = 12423.4995
$ws.Range("N125").Value = -22263.4995
$ws.Range("H138").Value = 11476.083
$ws.Range("I138").Value = 15965.714
$ws.Range("J138").Value = 5190.6
$ws.Range("K138").Value = 47897.142
$ws.Range("L138").Value = 15571.8
$ws.Range("M138").Value = -42757.142
$ws.Range("N138").Value = -25851.8
$ws.Range("N92").ClearContents()
$ws.Range("M94").ClearContents()
$ws.Range("N94").ClearContents()
$ws.Range("M97").ClearContents()
$ws.Range("N104").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2701
$ws.Range("I113").Value = 2944.1428
$ws.Range("J113").Value = 1850
$ws.Range("K113").Value = 2944.1428
$ws.Range("L113").Value = 1850
$ws.Range("M113").Value = -774.1428000000001
$ws.Range("N113").Value = -6190
$ws.Range("H122").Value = 6801.3
$ws.Range("I122").Value = 11402.6
$ws.Range("J122").Value = 2200
$ws.Range("K122").Value = 34207.8
$ws.Range("L122").Value = 6600
$ws.Range("M122").Value = -31757.8
$ws.Range("N122").Value = -11500

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H107").Value = 8000
$ws.Range("I107").Value = 8000
$ws.Range("K107").Value = 8000
$ws.Range("M107").Value = -6080
